# Apply scheduled-runner profit recalculation updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1266.625
$ws.Range("J46").Value = 1800
$ws.Range("L46").Value = 5400
$ws.Range("N46").Value = -5638
$ws.Range("H60").Value = 1266.625
$ws.Range("J60").Value = 1800
$ws.Range("L60").Value = 5400
$ws.Range("N60").Value = -6368
$ws.Range("H112").Value = 1257.01
$ws.Range("J112").Value = 1363.0449
$ws.Range("L112").Value = 4089.134700000001
$ws.Range("N112").Value = -6305.134700000001
$ws.Range("H137").Value = 5884770.5
$ws.Range("I137").Value = 8335416.5
$ws.Range("J137").Value = 3220
$ws.Range("K137").Value = 25006249.5
$ws.Range("L137").Value = 9660
$ws.Range("M137").Value = -25003699.5
$ws.Range("N137").Value = -14760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3035.18
$ws.Range("I32").Value = 2487.7092
$ws.Range("J32").Value = 6398.2144
$ws.Range("K32").Value = 2487.7092
$ws.Range("L32").Value = 6398.2144
$ws.Range("M32").Value = -2200.7092
$ws.Range("N32").Value = -6972.2144
$ws.Range("H61").Value = 2754.077
$ws.Range("I61").Value = 1753.25
$ws.Range("J61").Value = 3198.889
$ws.Range("K61").Value = 1753.25
$ws.Range("L61").Value = 3198.889
$ws.Range("M61").Value = -1541.25
$ws.Range("N61").Value = -3622.889
$ws.Range("H74").Value = 861.1429
$ws.Range("I74").Value = 800
$ws.Range("J74").Value = 907
$ws.Range("K74").Value = 800
$ws.Range("L74").Value = 907
$ws.Range("M74").Value = 74
$ws.Range("N74").Value = -2655
$ws.Range("H77").Value = 861.1429
$ws.Range("I77").Value = 800
$ws.Range("J77").Value = 907
$ws.Range("K77").Value = 4000
$ws.Range("L77").Value = 4535
$ws.Range("M77").Value = 368
$ws.Range("N77").Value = -13271
$ws.Range("H136").Value = 2754.077
$ws.Range("I136").Value = 1753.25
$ws.Range("J136").Value = 3198.889
$ws.Range("K136").Value = 5259.75
$ws.Range("L136").Value = 9596.667000000001
$ws.Range("M136").Value = -2709.75
$ws.Range("N136").Value = -14696.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3322.182
$ws.Range("I134").Value = 3120.4614
$ws.Range("K134").Value = 9361.3842
$ws.Range("M134").Value = -6826.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1820521
$ws.Range("I31").Value = 2129247.5
$ws.Range("J31").Value = 6753.125
$ws.Range("K31").Value = 2129247.5
$ws.Range("L31").Value = 6753.125
$ws.Range("M31").Value = -2128952.5
$ws.Range("N31").Value = -7343.125
$ws.Range("H34").Value = 1820521
$ws.Range("I34").Value = 2129247.5
$ws.Range("J34").Value = 6753.125
$ws.Range("K34").Value = 2129247.5
$ws.Range("L34").Value = 6753.125
$ws.Range("M34").Value = -2129045.5
$ws.Range("N34").Value = -7157.125
$ws.Range("H58").Value = 18521092
$ws.Range("I58").Value = 1049.1
$ws.Range("J58").Value = 29415236
$ws.Range("K58").Value = 1049.1
$ws.Range("L58").Value = 29415236
$ws.Range("M58").Value = -846.0999999999999
$ws.Range("N58").Value = -29415642
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H132").Value = 2697.1614
$ws.Range("I132").Value = 1780.5
$ws.Range("J132").Value = 5840
$ws.Range("K132").Value = 5341.5
$ws.Range("L132").Value = 17520
$ws.Range("M132").Value = -2811.5
$ws.Range("N132").Value = -22580
$ws.Range("H134").Value = 1966.5
$ws.Range("I134").Value = 1040.6
$ws.Range("J134").Value = 3123.875
$ws.Range("K134").Value = 3121.8
$ws.Range("L134").Value = 9371.625
$ws.Range("M134").Value = -586.7999999999997
$ws.Range("N134").Value = -14441.625
$ws.Range("H136").Value = 18521092
$ws.Range("I136").Value = 1049.1
$ws.Range("J136").Value = 29415236
$ws.Range("K136").Value = 3147.3
$ws.Range("L136").Value = 88245708
$ws.Range("M136").Value = -597.2999999999997
$ws.Range("N136").Value = -88250808

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5072.087
$ws.Range("I133").Value = 8757.5
$ws.Range("J133").Value = 4296.2104
$ws.Range("K133").Value = 26272.5
$ws.Range("L133").Value = 12888.6312
$ws.Range("M133").Value = -21212.5
$ws.Range("N133").Value = -23008.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2179.634
$ws.Range("I126").Value = 1430.8334
$ws.Range("J126").Value = 2308
$ws.Range("K126").Value = 4292.5002
$ws.Range("L126").Value = 6924
$ws.Range("M126").Value = -1822.5002
$ws.Range("N126").Value = -11864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2850.675
$ws.Range("I122").Value = 2493.276
$ws.Range("J122").Value = 3792.9092
$ws.Range("K122").Value = 7479.828
$ws.Range("L122").Value = 11378.7276
$ws.Range("M122").Value = -5029.828
$ws.Range("N122").Value = -16278.7276
$ws.Range("H132").Value = 2943.4333
$ws.Range("I132").Value = 1527.0667
$ws.Range("J132").Value = 4359.8
$ws.Range("K132").Value = 4581.2001
$ws.Range("L132").Value = 13079.4
$ws.Range("M132").Value = -2051.2001
$ws.Range("N132").Value = -18139.4
$ws.Range("H136").Value = 3127546.5
$ws.Range("I136").Value = 4764047
$ws.Range("K136").Value = 14292141
$ws.Range("M136").Value = -14289591

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16575
$ws.Range("J45").Value = 16575
$ws.Range("L45").Value = 16575
$ws.Range("N45").Value = -17557
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = $null
$ws.Range("H132").Value = 291954.56
$ws.Range("I132").Value = 387631.22
$ws.Range("K132").Value = 1162893.66
$ws.Range("M132").Value = -1160363.66
$ws.Range("H136").Value = 1803.8096
$ws.Range("I136").Value = 955.7143
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 2867.1429
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -317.1428999999998
$ws.Range("N136").Value = -15600
